# simplify steel description (remove RME)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The industrial mapping text in B2 lists construction-material shares;
# drop the stray "/RME" token from the steel (S) row so it reads
# "11.4% S/LFM+CDH/H:1" instead of "11.4% S/LFM+CDH/RME/H:1".
$cell = $ws.Range("B2")
$cell.Value = $cell.Value2 -replace "11\.4% S/LFM\+CDH/RME/H:1", "11.4% S/LFM+CDH/H:1"

# Wrap the long multi-line description and let the row grow to show it.
$cell.WrapText = $true
$ws.Rows.Item(2).RowHeight = 409.6

# Reflect the saved selection state (B2:B12, active cell at the bottom).
$ws.Range("B2:B12").Select()
